# Applies the numeric updates captured in the Hades_Profits.xlsx diff
# (crafting leve profit/price recalculations across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 832.6667
$ws.Range("I100").Value = 832.6667
$ws.Range("K100").Value = 832.6667
$ws.Range("M100").Value = -291.6667

$ws.Range("H112").Value = 1532.6595
$ws.Range("J112").Value = 1631.279
$ws.Range("L112").Value = 4893.837
$ws.Range("N112").Value = -7109.837

$ws.Range("H137").Value = 3228034
$ws.Range("I137").Value = 4763459.5
$ws.Range("J137").Value = 3640.3
$ws.Range("K137").Value = 14290378.5
$ws.Range("L137").Value = 10920.9
$ws.Range("M137").Value = -14287828.5
$ws.Range("N137").Value = -16020.9

$ws.Range("H141").Value = 3353.476
$ws.Range("I141").Value = 2338.9375
$ws.Range("J141").Value = 6600
$ws.Range("K141").Value = 7016.8125
$ws.Range("L141").Value = 19800
$ws.Range("M141").Value = -1836.8125
$ws.Range("N141").Value = -30160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5149908.5
$ws.Range("I32").Value = 6311449.5
$ws.Range("J32").Value = 19767.334
$ws.Range("K32").Value = 6311449.5
$ws.Range("L32").Value = 19767.334
$ws.Range("M32").Value = -6311162.5
$ws.Range("N32").Value = -20341.334

$ws.Range("H37").Value = 1869.5
$ws.Range("I37").Value = 1869.5
$ws.Range("K37").Value = 1869.5
$ws.Range("M37").Value = -1596.5

$ws.Range("H61").Value = 77078190
$ws.Range("I61").Value = 100101256
$ws.Range("J61").Value = 334673.34
$ws.Range("K61").Value = 100101256
$ws.Range("L61").Value = 334673.34
$ws.Range("M61").Value = -100101044
$ws.Range("N61").Value = -335097.34

$ws.Range("H74").Value = 10501215
$ws.Range("I74").Value = 14765618
$ws.Range("J74").Value = 144807.14
$ws.Range("K74").Value = 14765618
$ws.Range("L74").Value = 144807.14
$ws.Range("M74").Value = -14764744
$ws.Range("N74").Value = -146555.14

$ws.Range("H77").Value = 10501215
$ws.Range("I77").Value = 14765618
$ws.Range("J77").Value = 144807.14
$ws.Range("K77").Value = 73828090
$ws.Range("L77").Value = 724035.7000000001
$ws.Range("M77").Value = -73823722
$ws.Range("N77").Value = -732771.7000000001

$ws.Range("H132").Value = 54485.08
$ws.Range("I132").Value = 38650.63
$ws.Range("J132").Value = 93351.45
$ws.Range("K132").Value = 115951.89
$ws.Range("L132").Value = 280054.35
$ws.Range("M132").Value = -113421.89
$ws.Range("N132").Value = -285114.35

$ws.Range("H136").Value = 77078190
$ws.Range("I136").Value = 100101256
$ws.Range("J136").Value = 334673.34
$ws.Range("K136").Value = 300303768
$ws.Range("L136").Value = 1004020.02
$ws.Range("M136").Value = -300301218
$ws.Range("N136").Value = -1009120.02

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 31439068
$ws.Range("I105").Value = 35930136
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 35930136
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = -35928389
$ws.Range("N105").Value = -5094

$ws.Range("H134").Value = 3178.05
$ws.Range("I134").Value = 2671.8235
$ws.Range("J134").Value = 6046.6665
$ws.Range("K134").Value = 8015.470499999999
$ws.Range("L134").Value = 18139.9995
$ws.Range("M134").Value = -5480.470499999999
$ws.Range("N134").Value = -23209.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 52905.75
$ws.Range("I132").Value = 2508.3125
$ws.Range("K132").Value = 7524.9375
$ws.Range("M132").Value = -4994.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3426.3635
$ws.Range("J76").Value = 3469
$ws.Range("L76").Value = 10407
$ws.Range("N76").Value = -11173

$ws.Range("H79").Value = 3426.3635
$ws.Range("J79").Value = 3469
$ws.Range("L79").Value = 10407
$ws.Range("N79").Value = -13059

$ws.Range("H88").Value = 3683.3333
$ws.Range("J88").Value = 3683.3333
$ws.Range("L88").Value = 11049.9999
$ws.Range("N88").Value = -11905.9999

$ws.Range("H91").Value = 3683.3333
$ws.Range("J91").Value = 3683.3333
$ws.Range("L91").Value = 11049.9999
$ws.Range("N91").Value = -14013.9999

$ws.Range("H94").Value = 3298.4211
$ws.Range("I94").Value = 200
$ws.Range("J94").Value = 3470.5557
$ws.Range("K94").Value = 600
$ws.Range("L94").Value = 10411.6671
$ws.Range("M94").Value = 76
$ws.Range("N94").Value = -11763.6671

$ws.Range("H107").Value = 857.7846
$ws.Range("I107").Value = 393.70834
$ws.Range("J107").Value = 2168.1177
$ws.Range("K107").Value = 1181.12502
$ws.Range("L107").Value = 6504.353099999999
$ws.Range("M107").Value = 738.8749800000001
$ws.Range("N107").Value = -10344.3531

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0

$ws.Range("H132").Value = 2446.7896
$ws.Range("I132").Value = 2056.75
$ws.Range("J132").Value = 2730.4546
$ws.Range("K132").Value = 18510.75
$ws.Range("L132").Value = 24574.0914
$ws.Range("M132").Value = -15980.75
$ws.Range("N132").Value = -29634.0914

$ws.Range("H133").Value = 3252
$ws.Range("I133").Value = 3252
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 9756
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -4696

$ws.Range("H134").Value = 2414.1667
$ws.Range("I134").Value = 2088.182
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 6264.545999999999
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -1194.545999999999
$ws.Range("N134").Value = -28140

$ws.Range("H140").Value = 1949.7172
$ws.Range("I140").Value = 861.875
$ws.Range("J140").Value = 2159.4216
$ws.Range("K140").Value = 2585.625
$ws.Range("L140").Value = 6478.264800000001
$ws.Range("M140").Value = 2594.375
$ws.Range("N140").Value = -16838.2648

$ws.Range("H141").Value = 30000
$ws.Range("I141").Value = 30000
$ws.Range("K141").Value = 90000
$ws.Range("M141").Value = -84820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I93").Value = 2331.3333
$ws.Range("J93").Value = 2291.5
$ws.Range("K93").Value = 2331.3333
$ws.Range("L93").Value = 2291.5
$ws.Range("M93").Value = -1083.3333
$ws.Range("N93").Value = -4787.5

$ws.Range("H136").Value = 65290.676
$ws.Range("I136").Value = 31372.676
$ws.Range("J136").Value = 161391.67
$ws.Range("K136").Value = 94118.02799999999
$ws.Range("L136").Value = 484175.01
$ws.Range("M136").Value = -91568.02799999999
$ws.Range("N136").Value = -489275.01

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 34000
$ws.Range("J108").Value = 34000
$ws.Range("L108").Value = 34000
$ws.Range("N108").Value = -41680

# Cells that are fully removed from the row after recalculation
$wb.Worksheets.Item("CUL").Range("N123").ClearContents()
$wb.Worksheets.Item("CUL").Range("N133").ClearContents()
